$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update job image references from blog_1.jpg to the newly uploaded blog_4.jpg
$ws.Range("A2").Value = "blog_4.jpg"
$ws.Range("A4").Value = "blog_4.jpg"

# Restore single-cell selection on A4 (matches the saved selection state)
$ws.Range("A4").Select()
